$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 3) mirroring the structure of row 2.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 42606.571261574078

$ws.Range("B3").Value = 76
$ws.Range("C3").Value = 88
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 5904
$ws.Range("H3").Value = 2131
$ws.Range("I3").Value = 127
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Named"
